$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D (Price) and E (Volume 1h) columns with latest crypto data.
# D-column values are forced to Text format before assignment so that
# numeric-looking strings (e.g. "1.000", "0.9995") are kept verbatim
# instead of being auto-converted to numbers by Excel's type coercion.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.873.31"
$ws.Range("E2").Value = "  +1.48%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.728.25"
$ws.Range("E3").Value = "  +0.21%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  +0.24%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.35"
$ws.Range("E5").Value = "  -1.07%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.24%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4820"
$ws.Range("E7").Value = "  -1.61%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2594"
$ws.Range("E8").Value = "  -0.52%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06167"
$ws.Range("E9").Value = "  -0.58%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.726.36"
$ws.Range("E10").Value = "  -0.02%  "

# Row 11
$ws.Range("E11").Value = "  +2.15%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06894"
$ws.Range("E12").Value = "  -1.04%  "

# Row 13
$ws.Range("E13").Value = "  -0.24%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.451"
$ws.Range("E14").Value = "  -1.67%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.92"
$ws.Range("E15").Value = "  -0.30%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  +0.31%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.655.88"
$ws.Range("E17").Value = "  +0.68%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("E18").Value = "  +0.28%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007105"
$ws.Range("E19").Value = "  -1.01%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.34"
$ws.Range("E20").Value = "  +0.15%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.950.24"
$ws.Range("E21").Value = "  +0.08%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.392"
$ws.Range("E22").Value = "  -1.25%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.405"
$ws.Range("E23").Value = "  -1.38%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.049"
$ws.Range("E24").Value = "  -1.18%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.63"
$ws.Range("E25").Value = "  +1.42%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.17"
$ws.Range("E26").Value = "  -0.58%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.791"
$ws.Range("E27").Value = "  +2.47%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "106.47"
$ws.Range("E28").Value = "  -0.12%  "

# Row 29
$ws.Range("E29").Value = "  -1.95%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.938"
$ws.Range("E30").Value = "  +0.45%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07906"
$ws.Range("E31").Value = "  -1.29%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.655"
$ws.Range("E32").Value = "  +0.35%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04580"
$ws.Range("E33").Value = "  +2.12%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.595"
$ws.Range("E34").Value = "  -0.10%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9985"
$ws.Range("E35").Value = "  -0.44%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6152"
$ws.Range("E36").Value = "  -1.22%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9202"
$ws.Range("E37").Value = "  -2.00%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.468"
$ws.Range("E38").Value = "  +3.41%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.979"
$ws.Range("E39").Value = "  -0.84%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9995"
$ws.Range("E40").Value = "  +0.24%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.704"
$ws.Range("E41").Value = "  +5.58%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01492"
$ws.Range("E42").Value = "  +0.53%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.79"
$ws.Range("E43").Value = "  -0.07%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3820"
$ws.Range("E44").Value = "  -0.73%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.749"
$ws.Range("E45").Value = "  -2.10%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1149"
$ws.Range("E46").Value = "  -0.86%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05360"
$ws.Range("E47").Value = "  -0.14%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.833"
$ws.Range("E48").Value = "  +1.34%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.04"
$ws.Range("E49").Value = "  -1.78%  "

# Row 50
$ws.Range("E50").Value = "  +0.86%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.07"
$ws.Range("E51").Value = "  -0.79%  "
